$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$ws.Range("D4").Value = 0.701612903225806
$ws.Range("E4").Value = 0.683610867659947
$ws.Range("F4").Value = 0.728285077951002
$ws.Range("G4").Value = 0.568822553897181
$ws.Range("H4").Value = 0.726299694189602
$ws.Range("I4").Value = 0.667752442996743
$ws.Range("J4").Value = 0.65607476635514
$ws.Range("K4").Value = 0.706269349845201
$ws.Range("L4").Value = 0.685934489402698
$ws.Range("M4").Value = 0.668763102725367
$ws.Range("N4").Value = 0.602535832414553
$ws.Range("D5").Value = 0.578341013824885
$ws.Range("E5").Value = 0.609114811568799
$ws.Range("F5").Value = 0.682628062360802
$ws.Range("G5").Value = 0.681592039800995
$ws.Range("H5").Value = 0.697247706422018
$ws.Range("I5").Value = 0.51357220412595
$ws.Range("J5").Value = 0.536448598130841
$ws.Range("K5").Value = 0.498839009287926
$ws.Range("L5").Value = 0.765895953757225
$ws.Range("M5").Value = 0.590496156533892
$ws.Range("N5").Value = 0.42805953693495
$ws.Range("D6").Value = 0.668202764976959
$ws.Range("E6").Value = 0.670464504820333
$ws.Range("F6").Value = 0.736080178173719
$ws.Range("G6").Value = 0.648424543946932
$ws.Range("H6").Value = 0.723241590214067
$ws.Range("I6").Value = 0.642779587404995
$ws.Range("J6").Value = 0.598130841121495
$ws.Range("K6").Value = 0.506191950464396
$ws.Range("M6").Value = 0.560447239692523
$ws.Range("N6").Value = 0.612458654906284
$ws.Range("D7").Value = 0.432038834951456
$ws.Range("E7").Value = 0.385531135531136
$ws.Range("G7").Value = 0.287545787545788
$ws.Range("H7").Value = 0.437699680511182
$ws.Range("I7").Value = 0.394889663182346
$ws.Range("J7").Value = 0.323353293413174
$ws.Range("K7").Value = 0.216003316749585
$ws.Range("N7").Value = 0.372829728906488
$ws.Range("D8").Value = 0.327669902912621
$ws.Range("E8").Value = 0.335164835164835
$ws.Range("H8").Value = 0.400958466453674
$ws.Range("I8").Value = 0.35075493612079
$ws.Range("N8").Value = 0.30855924459336
$ws.Range("D9").Value = 0.327669902912621
$ws.Range("E9").Value = 0.451465201465201
$ws.Range("G9").Value = 0.41025641025641
$ws.Range("H9").Value = 0.501597444089457
$ws.Range("I9").Value = 0.451800232288037
$ws.Range("J9").Value = 0.389221556886228
$ws.Range("K9").Value = 0.313432835820896
$ws.Range("L9").Value = 0.39980732177264
$ws.Range("N9").Value = 0.411513859275053
